# "Generate Report for Handback"
#
# The qimu localization-status report records, per locale sheet, the
# timestamp of the most recently received handback file for each source
# document (column K = "Latest Handback DateTime"). This run recorded a
# fresh handback for the "61aa49de-59cc-4317-a7ed-0b964375a556" document
# (report row 2) in both locale sheets; the "2a9b8297-..." document
# (report row 3) did not receive a new handback in this pass, so its
# datetime is left untouched.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Latest Handback DateTime (column K, row 2) for the zh-cn target.
$zhcn.Range("K2").Value = "2016-10-10 09:54:47"

# Latest Handback DateTime (column K, row 2) for the de-de target.
$dede.Range("K2").Value = "2016-10-10 09:55:02"
